$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("liste")
$ws2 = $wb.Worksheets.Item("altri_dati")

# ---------------------------------------------------------------------------
# New minor-party rows (14-35) replacing the old placeholder "Altri 1/2/3"
# rows. Each row repeats the party name in A and B, is flagged "false" in E
# and tagged "ALTRI" in G, same pattern as the pre-existing rows 14/15.
# ---------------------------------------------------------------------------
$names = @(
    "ITALIA SOVRANA E POPOLARE",
    "VITA",
    "MASTELLA NOI DI CENTRO EUROPEISTI",
    "FREE",
    "ALTERNATIVA PER L'ITALIA - NO GREEN PASS",
    "VALLE D'AOSTA APERTA",
    "LEGA PER SALVINI PREMIER - FORZA ITALIA - NOI MODERATI - FRATELLI D'ITALIA",
    "LA RENAISSANCE VALDÔTAINE",
    "PARTITO COMUNISTA ITALIANO",
    "VALLÉE D’AOSTE – AUTONOMIE PROGRÈS FÉDÉRALISME",
    "SUD CHIAMA NORD",
    "PARTITO ANIMALISTA - UCDL - 10 VOLTE MEGLIO",
    "PARTITO DELLA FOLLIA CREATIVA",
    "FORZA DEL POPOLO",
    "SÜDTIROLER VOLKSPARTEI (SVP) - PATT",
    "PARTITO COMUNISTA DEI LAVORATORI",
    "DESTRE UNITE",
    "POUR L’AUTONOMIE – PER L’AUTONOMIA",
    "CAMPOBASE - +EUROPA - ALLEANZA VERDI E SINISTRA – PD - AZIONE-ITALIA VIVA",
    "PARTITO DEMOCRATICO - +EUROPA - ALLEANZA VERDI E SINISTRA",
    "DIE FREIHEITLICHEN",
    "TEAM K"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 14 + $i
    $ws1.Cells.Item($r, 1).Value = $names[$i]
    $ws1.Cells.Item($r, 2).Value = $names[$i]
    $ws1.Cells.Item($r, 5).Value = $false
    $ws1.Cells.Item($r, 7).Value = "ALTRI"
}

# C14 keeps its own (non-shared) formula; C15:C35 is filled as one block so
# the engine groups it into a single shared formula, matching the source.
$ws1.Range("C14").Formula = "=0.02/22"
$ws1.Range("C15:C35").Formula = "=0.02/22"
$ws1.Range("C14:C35").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Column widths (A widened to fit the longer party names, B no longer
# auto-best-fit).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 73.1
$ws1.Columns.Item(2).ColumnWidth = 15.67

# ---------------------------------------------------------------------------
# View state: "liste" becomes the active sheet/tab, selection moves to A28;
# "altri_dati" loses its tab-selected flag (selection on that sheet is
# untouched).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A28").Select() | Out-Null
